$wb = $excel.ActiveWorkbook

# --- logo sheet ---
$ws1 = $wb.Worksheets.Item("logo")
$ws1.Range("B2").Value = "yes"
$ws1.Range("B3").Select()

# --- firstbanner sheet ---
$ws2 = $wb.Worksheets.Item("firstbanner")
$ws2.Range("B2").Value = "yes"
$ws2.Range("B3").Select()

# --- trending sheet ---
$ws3 = $wb.Worksheets.Item("trending")
$ws3.Range("B2").Value = "yes"
$ws3.Range("B3").Select()

# --- introduction sheet ---
$ws4 = $wb.Worksheets.Item("introduction")
$ws4.Range("B2").Value = "yes"
$ws4.Range("B2").Select()

# --- levelUp sheet (footer page) ---
$ws5 = $wb.Worksheets.Item("levelUp")
$ws5.Range("E5").Value = "products/hydra-filling-cream-2"
$ws5.Range("E5").Select()

# Restore the originally active sheet/tab (levelUp was tabSelected in the source file)
$ws5.Activate()
